$wb = $excel.ActiveWorkbook

# --- "New Horizons" sheet: append rows 59-64 ---
$ws = $wb.Worksheets.Item("New Horizons")

# Row 59
$ws.Range("A57:N57").Copy() | Out-Null
$ws.Range("A59:N59").PasteSpecial(-4122) | Out-Null
$ws.Range("A59").Value2 = 'GENE'
$ws.Range("B59").Value2 = 'LILY'
$ws.Range("C59").Value2 = 'KIT'
$ws.Range("D59").Value2 = 'OLLIE'
$ws.Range("E59").Value2 = 'MEG'
$ws.Range("F59").Value2 = 'CORDELIUS'
$ws.Range("G59").Value2 = 'Equipo 2'
$ws.Range("H59").Value2 = 'SK|Ope'
$ws.Range("I59").Value2 = 'SK|Joker'
$ws.Range("J59").Value2 = 'SK|Yoshi825'
$ws.Range("K59").Value2 = 'HMB|BosS'
$ws.Range("L59").Value2 = 'HMB|Lukii'
$ws.Range("M59").Value2 = 'HMB|Symantec'
$ws.Range("N59").Value2 = '20250724T172257.000Z'

# Row 60
$ws.Range("A57:N57").Copy() | Out-Null
$ws.Range("A60:N60").PasteSpecial(-4122) | Out-Null
$ws.Range("A60").Value2 = 'BONNIE'
$ws.Range("B60").Value2 = 'LUMI'
$ws.Range("C60").Value2 = 'BROCK'
$ws.Range("D60").Value2 = 'ANGELO'
$ws.Range("E60").Value2 = 'GENE'
$ws.Range("F60").Value2 = 'HANK'
$ws.Range("G60").Value2 = 'Equipo 2'
$ws.Range("H60").Value2 = 'FUT|Nowy297'
$ws.Range("I60").Value2 = 'FUT|MeOw'
$ws.Range("J60").Value2 = 'FUT|GeRo'
$ws.Range("K60").Value2 = 'TH|iKaoss'
$ws.Range("L60").Value2 = 'TH|Zhar'
$ws.Range("M60").Value2 = 'TH|LeNain'
$ws.Range("N60").Value2 = '20250724T174707.000Z'

# Row 61
$ws.Range("A58:N58").Copy() | Out-Null
$ws.Range("A61:N61").PasteSpecial(-4122) | Out-Null
$ws.Range("A61").Value2 = 'BONNIE'
$ws.Range("B61").Value2 = 'LUMI'
$ws.Range("C61").Value2 = 'BROCK'
$ws.Range("D61").Value2 = 'ANGELO'
$ws.Range("E61").Value2 = 'GENE'
$ws.Range("F61").Value2 = 'HANK'
$ws.Range("G61").Value2 = 'Equipo 1'
$ws.Range("H61").Value2 = 'FUT|Nowy297'
$ws.Range("I61").Value2 = 'FUT|MeOw'
$ws.Range("J61").Value2 = 'FUT|GeRo'
$ws.Range("K61").Value2 = 'TH|iKaoss'
$ws.Range("L61").Value2 = 'TH|Zhar'
$ws.Range("M61").Value2 = 'TH|LeNain'
$ws.Range("N61").Value2 = '20250724T174359.000Z'

# Row 62
$ws.Range("A57:N57").Copy() | Out-Null
$ws.Range("A62:N62").PasteSpecial(-4122) | Out-Null
$ws.Range("A62").Value2 = 'BONNIE'
$ws.Range("B62").Value2 = 'LUMI'
$ws.Range("C62").Value2 = 'BROCK'
$ws.Range("D62").Value2 = 'ANGELO'
$ws.Range("E62").Value2 = 'GENE'
$ws.Range("F62").Value2 = 'HANK'
$ws.Range("G62").Value2 = 'Equipo 2'
$ws.Range("H62").Value2 = 'FUT|Nowy297'
$ws.Range("I62").Value2 = 'FUT|MeOw'
$ws.Range("J62").Value2 = 'FUT|GeRo'
$ws.Range("K62").Value2 = 'TH|iKaoss'
$ws.Range("L62").Value2 = 'TH|Zhar'
$ws.Range("M62").Value2 = 'TH|LeNain'
$ws.Range("N62").Value2 = '20250724T174137.000Z'

# Row 63
$ws.Range("A58:N58").Copy() | Out-Null
$ws.Range("A63:N63").PasteSpecial(-4122) | Out-Null
$ws.Range("A63").Value2 = 'GUS'
$ws.Range("B63").Value2 = 'BUSTER'
$ws.Range("C63").Value2 = 'KAZE'
$ws.Range("D63").Value2 = 'R-T'
$ws.Range("E63").Value2 = 'SQUEAK'
$ws.Range("F63").Value2 = 'ALLI'
$ws.Range("G63").Value2 = 'Equipo 1'
$ws.Range("H63").Value2 = 'FUT|GeRo'
$ws.Range("I63").Value2 = 'FUT|Nowy297'
$ws.Range("J63").Value2 = 'FUT|MeOw'
$ws.Range("K63").Value2 = 'TH|Zhar'
$ws.Range("L63").Value2 = 'TH|iKaoss'
$ws.Range("M63").Value2 = 'TH|LeNain'
$ws.Range("N63").Value2 = '20250724T173429.000Z'

# Row 64
$ws.Range("A58:N58").Copy() | Out-Null
$ws.Range("A64:N64").PasteSpecial(-4122) | Out-Null
$ws.Range("A64").Value2 = 'GUS'
$ws.Range("B64").Value2 = 'BUSTER'
$ws.Range("C64").Value2 = 'KAZE'
$ws.Range("D64").Value2 = 'R-T'
$ws.Range("E64").Value2 = 'SQUEAK'
$ws.Range("F64").Value2 = 'ALLI'
$ws.Range("G64").Value2 = 'Equipo 1'
$ws.Range("H64").Value2 = 'FUT|GeRo'
$ws.Range("I64").Value2 = 'FUT|Nowy297'
$ws.Range("J64").Value2 = 'FUT|MeOw'
$ws.Range("K64").Value2 = 'TH|Zhar'
$ws.Range("L64").Value2 = 'TH|iKaoss'
$ws.Range("M64").Value2 = 'TH|LeNain'
$ws.Range("N64").Value2 = '20250724T173209.000Z'

$excel.CutCopyMode = $false

# --- "Hot Potato" sheet: append row 73 ---
$ws2 = $wb.Worksheets.Item("Hot Potato")

# Row 73
$ws2.Range("A72:N72").Copy() | Out-Null
$ws2.Range("A73:N73").PasteSpecial(-4122) | Out-Null
$ws2.Range("A73").Value2 = 'CHARLIE'
$ws2.Range("B73").Value2 = 'LILY'
$ws2.Range("C73").Value2 = 'SHADE'
$ws2.Range("D73").Value2 = 'JESSIE'
$ws2.Range("E73").Value2 = 'MELODIE'
$ws2.Range("F73").Value2 = 'KAZE'
$ws2.Range("G73").Value2 = 'Equipo 1'
$ws2.Range("H73").Value2 = 'FUT|GeRo'
$ws2.Range("I73").Value2 = 'FUT|Nowy297'
$ws2.Range("J73").Value2 = 'FUT|MeOw'
$ws2.Range("K73").Value2 = 'TH|iKaoss'
$ws2.Range("L73").Value2 = 'TH|LeNain'
$ws2.Range("M73").Value2 = 'TH|Zhar'
$ws2.Range("N73").Value2 = '20250724T172343.000Z'

$excel.CutCopyMode = $false